$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # BOPIS_TestData
$ws2 = $wb.Worksheets.Item(2)   # SFS_TestData
$ws3 = $wb.Worksheets.Item(3)   # STS_TestData
$ws4 = $wb.Worksheets.Item(4)   # BackUp

# --- Data edits -----------------------------------------------------------

# BOPIS_TestData
$ws1.Range("C2").Value = 1

# SFS_TestData
$ws2.Range("B2").Value = 10817056
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 1

# STS_TestData
$ws3.Range("B2").Value = 91327491
$ws3.Range("C2").Value = 1
$ws3.Range("E2").Value = 83
$ws3.Range("F2").Value = 89109

# --- Selection / active sheet changes -------------------------------------

[void]$ws1.Range("D13").Select()
[void]$ws2.Range("B2").Select()
[void]$ws3.Range("B2").Select()
[void]$ws3.Activate()
